# Apply new TPM data: remove Inflammatory-Mac target-cluster rows and recompute values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows for the removed "Inflammatory-Mac" target cluster (bottom-up)
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(10).Delete()

# Row 2
$ws.Range("I2").Value = 0.1477750351608889
$ws.Range("J2").Value = 0.1477750351608889
$ws.Range("M2").Value = 2.544438666666667
$ws.Range("N2").Value = 7.633316000000001
$ws.Range("O2").Value = 0.201325300207035
$ws.Range("P2").Value = 0.201325300207035
$ws.Range("Q2").Value = 0.1855379231346667
$ws.Range("R2").Value = 1.669841308212
$ws.Range("S2").Value = 0.02975085331687111
$ws.Range("T2").Value = 0.02975085331687111

# Row 3
$ws.Range("I3").Value = 0.1477750351608889
$ws.Range("J3").Value = 0.1477750351608889
$ws.Range("O3").Value = 0.6969390273602759
$ws.Range("P3").Value = 0.696939027360276
$ws.Range("S3").Value = 0.1029901892731605
$ws.Range("T3").Value = 0.1029901892731605

# Row 4
$ws.Range("D4").Value = "MuSCs"
$ws.Range("I4").Value = 0.1477750351608889
$ws.Range("J4").Value = 0.1477750351608889
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.273916333333333
$ws.Range("N4").Value = 3.821749
$ws.Range("O4").Value = 0.1007969229547075
$ws.Range("P4").Value = 0.1007969229547075
$ws.Range("Q4").Value = 0.09289270511033333
$ws.Range("R4").Value = 0.8360343459930001
$ws.Range("S4").Value = 0.01489526883374131
$ws.Range("T4").Value = 0.01489526883374131

# Row 5
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("I5").Value = 0.1477750351608889
$ws.Range("J5").Value = 0.1477750351608889
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01186433333333333
$ws.Range("N5").Value = 0.035593
$ws.Range("O5").Value = 0.0009387494779816524
$ws.Range("P5").Value = 0.0009387494779816526
$ws.Range("Q5").Value = 0.0008651353223333333
$ws.Range("R5").Value = 0.007786217901
$ws.Range("S5").Value = 0.0001387237371160048
$ws.Range("T5").Value = 0.0001387237371160048

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.420527
$ws.Range("H6").Value = 1.261581
$ws.Range("I6").Value = 0.852224964839111
$ws.Range("J6").Value = 0.852224964839111
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.544438666666667
$ws.Range("N6").Value = 7.633316000000001
$ws.Range("O6").Value = 0.201325300207035
$ws.Range("P6").Value = 0.201325300207035
$ws.Range("Q6").Value = 1.070005159177333
$ws.Range("R6").Value = 9.630046432596002
$ws.Range("S6").Value = 0.1715744468901639
$ws.Range("T6").Value = 0.1715744468901639

# Row 7
$ws.Range("D7").Value = "FAPs"
$ws.Range("G7").Value = 0.420527
$ws.Range("H7").Value = 1.261581
$ws.Range("I7").Value = 0.852224964839111
$ws.Range("J7").Value = 0.852224964839111
$ws.Range("M7").Value = 8.808225333333333
$ws.Range("N7").Value = 26.424676
$ws.Range("O7").Value = 0.6969390273602759
$ws.Range("P7").Value = 0.696939027360276
$ws.Range("Q7").Value = 3.704096574750667
$ws.Range("R7").Value = 33.336869172756
$ws.Range("S7").Value = 0.5939488380871154
$ws.Range("T7").Value = 0.5939488380871155

# Row 8
$ws.Range("D8").Value = "MuSCs"
$ws.Range("G8").Value = 0.420527
$ws.Range("H8").Value = 1.261581
$ws.Range("I8").Value = 0.852224964839111
$ws.Range("J8").Value = 0.852224964839111
$ws.Range("M8").Value = 1.273916333333333
$ws.Range("N8").Value = 3.821749
$ws.Range("O8").Value = 0.1007969229547075
$ws.Range("P8").Value = 0.1007969229547075
$ws.Range("Q8").Value = 0.5357162139076668
$ws.Range("R8").Value = 4.821445925169001
$ws.Range("S8").Value = 0.08590165412096615
$ws.Range("T8").Value = 0.08590165412096616

# Row 9
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("G9").Value = 0.420527
$ws.Range("H9").Value = 1.261581
$ws.Range("I9").Value = 0.852224964839111
$ws.Range("J9").Value = 0.852224964839111
$ws.Range("M9").Value = 0.01186433333333333
$ws.Range("N9").Value = 0.035593
$ws.Range("O9").Value = 0.0009387494779816524
$ws.Range("P9").Value = 0.0009387494779816526
$ws.Range("Q9").Value = 0.004989272503666667
$ws.Range("R9").Value = 0.044903452533
$ws.Range("S9").Value = 0.0008000257408656476
$ws.Range("T9").Value = 0.0008000257408656478

